$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4128389656543732
$ws.Range("B1").Value = 0.6896315217018127
$ws.Range("C1").Value = 0.8533981442451477
$ws.Range("D1").Value = 4.240665912628174
$ws.Range("E1").Value = 1.286994814872742
